$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds date-like text ("YYYY-MM-DD"); pre-format as Text so
# Excel does not silently reinterpret it as a date serial number.
$ws.Range("C25:C30").NumberFormat = "@"

# Row 25
$ws.Range("A25").Value = "WC47 NACP"
$ws.Range("B25").Value = "Fallo tolva"
$ws.Range("C25").Value = "2024-06-05"
$ws.Range("D25").Value = "10:48:47"
$ws.Range("E25").Value = "Mañana"
$ws.Range("F25").Value = "10:48:48"
$ws.Range("G25").Value = "0:00:01"
$ws.Range("H25").Value = "N/A"

# Row 26
$ws.Range("A26").Value = "WC47 NACP"
$ws.Range("B26").Value = "No coge placa"
$ws.Range("C26").Value = "2024-06-05"
$ws.Range("D26").Value = "10:48:49"
$ws.Range("E26").Value = "Mañana"
$ws.Range("F26").Value = "10:48:50"
$ws.Range("G26").Value = "0:00:01"
$ws.Range("H26").Value = "N/A"

# Row 27
$ws.Range("A27").Value = "WC47 NACP"
$ws.Range("B27").Value = "Fallo en elevador"
$ws.Range("C27").Value = "2024-06-05"
$ws.Range("D27").Value = "10:48:52"
$ws.Range("E27").Value = "Mañana"
$ws.Range("F27").Value = "10:48:52"
$ws.Range("G27").Value = "0:00:00"
$ws.Range("H27").Value = "N/A"

# Row 28
$ws.Range("A28").Value = "WC48 P5F"
$ws.Range("B28").Value = "AOI (fallo etiqueta)"
$ws.Range("C28").Value = "2024-06-05"
$ws.Range("D28").Value = "10:50:52"
$ws.Range("E28").Value = "Mañana"
$ws.Range("F28").Value = "10:50:53"
$ws.Range("G28").Value = "0:00:01"
$ws.Range("H28").Value = "N/A"

# Row 29
$ws.Range("A29").Value = "WC48 P5F"
$ws.Range("B29").Value = "Cámara no detecta Pcb"
$ws.Range("C29").Value = "2024-06-05"
$ws.Range("D29").Value = "10:50:54"
$ws.Range("E29").Value = "Mañana"
$ws.Range("F29").Value = "10:50:57"
$ws.Range("G29").Value = "0:00:03"
$ws.Range("H29").Value = "N/A"

# Row 30
$ws.Range("A30").Value = "WC48 P5F"
$ws.Range("B30").Value = "No detecta presencia power CP"
$ws.Range("C30").Value = "2024-06-05"
$ws.Range("D30").Value = "10:50:56"
$ws.Range("E30").Value = "Mañana"
$ws.Range("F30").Value = "10:50:57"
$ws.Range("G30").Value = "0:00:01"
$ws.Range("H30").Value = "N/A"
